$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter new place data row (PlaceID 1004) directly into row 10,
# leaving rows 5-9 empty (sparse), matching how the data was typed in Excel.
$ws.Range("A10").Value = 1004

# Move the active selection to B10, as if the user pressed Tab/Right after
# typing the value in A10.
$ws.Range("B10").Select()
